$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item(1)
$wsZhCn = $wb.Worksheets.Item(2)
$wsDeDe = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    This text is shared by every row/column that shows the handoff status,
#    so every cell currently holding it needs to be rewritten.
# ---------------------------------------------------------------------------
$handedBack = "Handed back: in sync with en-US"
$wsOverview.Range("E2").Value = $handedBack
$wsOverview.Range("F2").Value = $handedBack
$wsOverview.Range("E3").Value = $handedBack
$wsOverview.Range("F3").Value = $handedBack
$wsZhCn.Range("C2").Value = $handedBack
$wsZhCn.Range("C3").Value = $handedBack
$wsDeDe.Range("C2").Value = $handedBack
$wsDeDe.Range("C3").Value = $handedBack

# ---------------------------------------------------------------------------
# 2. zh-cn sheet: fill in the "Latest Target File" (I) and
#    "Latest Handback File" (J) columns now that the handback finished.
# ---------------------------------------------------------------------------
$wsZhCn.Range("J2").Value = "675366d8-4363-458e-836e-4ad4bab87bb1.13e578bafcbbff652befacf48c3edd2d64cad9f9.zh-cn.xlf"
$wsZhCn.Range("J3").Value = "b13875ee-b2af-4acc-bc4a-cc938995b2ed.858577ad3584e13f1c01f698e93983cf7da44ed4.zh-cn.xlf"

# ---------------------------------------------------------------------------
# 3. de-de sheet: same Latest Target File / Latest Handback File update, plus
#    the Latest Handback DateTime (K) now has a real timestamp.
# ---------------------------------------------------------------------------
$wsDeDe.Range("J2").Value = "675366d8-4363-458e-836e-4ad4bab87bb1.13e578bafcbbff652befacf48c3edd2d64cad9f9.de-de.xlf"
$wsDeDe.Range("J3").Value = "b13875ee-b2af-4acc-bc4a-cc938995b2ed.858577ad3584e13f1c01f698e93983cf7da44ed4.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-31 13:35:40"
$wsDeDe.Range("K3").Value = "2016-08-31 13:35:40"

# The zh-cn handback timestamp is recorded a little earlier than de-de's.
$wsZhCn.Range("K2").Value = "2016-08-31 13:35:23"
$wsZhCn.Range("K3").Value = "2016-08-31 13:35:23"

# ---------------------------------------------------------------------------
# 4. Re-create the hyperlinks on both language sheets: keep the existing
#    "Source File Name" (A) links, and add matching links on the newly
#    populated "Latest Target File" (I) column pointing at the same source
#    markdown file on GitHub.
# ---------------------------------------------------------------------------
foreach ($ws in @($wsZhCn, $wsDeDe)) {
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1f54c914667879b9de2efc580650aa8cc67d9b5f/e2e/675366d8-4363-458e-836e-4ad4bab87bb1.md", "", "", "675366d8-4363-458e-836e-4ad4bab87bb1.md")
    $ws.Hyperlinks.Add($ws.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1f54c914667879b9de2efc580650aa8cc67d9b5f/e2e/675366d8-4363-458e-836e-4ad4bab87bb1.md", "", "", "675366d8-4363-458e-836e-4ad4bab87bb1.md")
    $ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1f54c914667879b9de2efc580650aa8cc67d9b5f/e2e/b13875ee-b2af-4acc-bc4a-cc938995b2ed.md", "", "", "b13875ee-b2af-4acc-bc4a-cc938995b2ed.md")
    $ws.Hyperlinks.Add($ws.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1f54c914667879b9de2efc580650aa8cc67d9b5f/e2e/b13875ee-b2af-4acc-bc4a-cc938995b2ed.md", "", "", "b13875ee-b2af-4acc-bc4a-cc938995b2ed.md")
}

# ---------------------------------------------------------------------------
# 5. Widen the columns that now show full file names, to match the report
#    generator's updated layout.
# ---------------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 29.09
$wsOverview.Columns.Item(6).ColumnWidth = 29.09

foreach ($ws in @($wsZhCn, $wsDeDe)) {
    $ws.Columns.Item(3).ColumnWidth = 29.09
    $ws.Columns.Item(9).ColumnWidth = 39.17
    $ws.Columns.Item(10).ColumnWidth = 39.17
}
